# Fruta / hortaliza, semanal
# Insert a new weekly price record before current row 20 (Provincia del Elquí,
# 2022-08-26). Excel's native row insert shifts every existing row 20..43
# down to 21..44, which reproduces the rest of the diff automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(20).Insert()

$ws.Range("A20").Value = 7
$ws.Range("B20").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C20").Value = 'Ñuble'
$ws.Range("D20").Value = 44799
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112026
$ws.Range("G20").Value = 'Haba'
$ws.Range("H20").Value = 'Sin especificar'
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 11000
$ws.Range("L20").Value = 11000
$ws.Range("M20").Value = 11000
$ws.Range("N20").Value = '$/saco 25 kilos'
$ws.Range("O20").Value = 'Provincia del Elquí'
$ws.Range("P20").Value = 440
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = 'Hortaliza'

Write-Host "Inserted new row 20; sheet now spans" $ws.UsedRange.Address()
